$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB6").Value = 23
$ws.Range("AD6").Value = 9.5
$ws.Range("AM6").Value = 51
$ws.Range("AN6").Value = 3.6
$ws.Range("AQ6").Value = 21
$ws.Range("AU6").Value = 8
$ws.Range("AX6").Value = 29
$ws.Range("BA6").Value = 101
$ws.Range("G6").Value = 1.45
$ws.Range("H6").Value = 4.33
$ws.Range("I6").Value = 6.25
$ws.Range("J6").Value = 1.95
$ws.Range("O6").Value = 1.2
$ws.Range("P6").Value = 4.5
$ws.Range("Q6").Value = 1.62
$ws.Range("R6").Value = 2.3
$ws.Range("W6").Value = 9
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 12
$ws.Range("Q7").Value = 1.53
$ws.Range("R7").Value = 2.5
$ws.Range("AM9").Value = 23
$ws.Range("G9").Value = 3.3
$ws.Range("I9").Value = 2.2
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 2.02
$ws.Range("AE10").Value = 21
$ws.Range("AF10").Value = 81
$ws.Range("AH10").Value = 10
$ws.Range("AI10").Value = 23
$ws.Range("AJ10").Value = 17
$ws.Range("AO10").Value = 10
$ws.Range("AQ10").Value = 34
$ws.Range("AU10").Value = 9.5
$ws.Range("AV10").Value = 81
$ws.Range("AW10").Value = 6.5
$ws.Range("AX10").Value = 29
$ws.Range("BD10").Value = 151
$ws.Range("G10").Value = 1.8
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 2.5
$ws.Range("L10").Value = 5.5
$ws.Range("U10").Value = 2.2
$ws.Range("V10").Value = 1.62
$ws.Range("AA11").Value = 15
$ws.Range("AB11").Value = 41
$ws.Range("AC11").Value = 8.5
$ws.Range("AD11").Value = 9
$ws.Range("AE11").Value = 29
$ws.Range("AF11").Value = 101
$ws.Range("AH11").Value = 15
$ws.Range("AJ11").Value = 26
$ws.Range("AN11").Value = 3.1
$ws.Range("AP11").Value = 23
$ws.Range("AS11").Value = 251
$ws.Range("AT11").Value = 2.63
$ws.Range("AU11").Value = 11
$ws.Range("AW11").Value = 9
$ws.Range("AY11").Value = 51
$ws.Range("AZ11").Value = 251
$ws.Range("BA11").Value = 301
$ws.Range("BD11").Value = 126
$ws.Range("G11").Value = 1.4
$ws.Range("H11").Value = 4.33
$ws.Range("I11").Value = 8.5
$ws.Range("K11").Value = 2.2
$ws.Range("L11").Value = 8.5
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 1.36
$ws.Range("P11").Value = 3
$ws.Range("Q11").Value = 2.1
$ws.Range("R11").Value = 1.7
$ws.Range("S11").Value = 1.44
$ws.Range("T11").Value = 2.63
$ws.Range("U11").Value = 2.5
$ws.Range("V11").Value = 1.5
$ws.Range("W11").Value = 5
$ws.Range("X11").Value = 5.5
$ws.Range("Y11").Value = 9
$ws.Range("Z11").Value = 8.5
$ws.Range("AF13").Value = 41
$ws.Range("AG13").Value = 151
$ws.Range("AI13").Value = 23
$ws.Range("AM13").Value = 34
$ws.Range("AP13").Value = 17
$ws.Range("AS13").Value = 101
$ws.Range("AX13").Value = 23
$ws.Range("BA13").Value = 81
$ws.Range("BB13").Value = 151
$ws.Range("G13").Value = 1.7
$ws.Range("I13").Value = 4.5
$ws.Range("L13").Value = 4.75
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 15
$ws.Range("O13").Value = 1.2
$ws.Range("P13").Value = 4.33
$ws.Range("Q13").Value = 1.7
$ws.Range("R13").Value = 2.1
$ws.Range("U13").Value = 1.67
$ws.Range("V13").Value = 2.1
$ws.Range("W13").Value = 8.5
$ws.Range("X13").Value = 9
$ws.Range("AB14").Value = 26
$ws.Range("AD14").Value = 29
$ws.Range("AF14").Value = 51
$ws.Range("AG14").Value = 151
$ws.Range("AH14").Value = 67
$ws.Range("AI14").Value = 126
$ws.Range("AL14").Value = 126
$ws.Range("AM14").Value = 67
$ws.Range("G14").Value = 1.06
$ws.Range("H14").Value = 12
$ws.Range("I14").Value = 23
$ws.Range("Q14").Value = 1.11
$ws.Range("R14").Value = 6.5
$ws.Range("W14").Value = 26
$ws.Range("X14").Value = 13
$ws.Range("Z14").Value = 10
$ws.Range("AA15").Value = 29
$ws.Range("AB15").Value = 41
$ws.Range("AH15").Value = 7
$ws.Range("AI15").Value = 10
$ws.Range("AK15").Value = 19
$ws.Range("AO15").Value = 19
$ws.Range("AQ15").Value = 67
$ws.Range("AW15").Value = 4
$ws.Range("AX15").Value = 12
$ws.Range("BB15").Value = 151
$ws.Range("G15").Value = 3.3
$ws.Range("I15").Value = 2.15
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 2.1
$ws.Range("L15").Value = 2.88
$ws.Range("W15").Value = 9.5
$ws.Range("X15").Value = 17
$ws.Range("AA16").Value = 29
$ws.Range("AC16").Value = 6
$ws.Range("AD16").Value = 5.5
$ws.Range("AG16").Value = 501
$ws.Range("AH16").Value = 6.5
$ws.Range("AJ16").Value = 11
$ws.Range("AK16").Value = 26
$ws.Range("AL16").Value = 26
$ws.Range("AM16").Value = 41
$ws.Range("AO16").Value = 19
$ws.Range("AP16").Value = 34
$ws.Range("AQ16").Value = 67
$ws.Range("AS16").Value = 351
$ws.Range("AT16").Value = 2.25
$ws.Range("AU16").Value = 9
$ws.Range("AV16").Value = 81
$ws.Range("AW16").Value = 4.5
$ws.Range("AX16").Value = 17
$ws.Range("BA16").Value = 101
$ws.Range("BB16").Value = 301
$ws.Range("H16").Value = 2.88
$ws.Range("I16").Value = 2.7
$ws.Range("K16").Value = 1.91
$ws.Range("L16").Value = 3.5
$ws.Range("M16").Value = 1.13
$ws.Range("N16").Value = 6
$ws.Range("O16").Value = 1.53
$ws.Range("P16").Value = 2.5
$ws.Range("Q16").Value = 2.7
$ws.Range("R16").Value = 1.44
$ws.Range("S16").Value = 1.57
$ws.Range("T16").Value = 2.25
$ws.Range("U16").Value = 2.1
$ws.Range("V16").Value = 1.67
$ws.Range("W16").Value = 7
$ws.Range("Y16").Value = 12
$ws.Range("Z16").Value = 34
$ws.Range("AA17").Value = 23
$ws.Range("AE17").Value = 13
$ws.Range("AG17").Value = 151
$ws.Range("AI17").Value = 11
$ws.Range("AK17").Value = 21
$ws.Range("AL17").Value = 17
$ws.Range("AO17").Value = 17
$ws.Range("AP17").Value = 23
$ws.Range("AR17").Value = 67
$ws.Range("AX17").Value = 12
$ws.Range("G17").Value = 3.2
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 2.15
$ws.Range("Q17").Value = 1.86
$ws.Range("R17").Value = 2.04
$ws.Range("U17").Value = 1.67
$ws.Range("V17").Value = 2.1
$ws.Range("Y17").Value = 11
$ws.Range("AC18").Value = 10
$ws.Range("AE18").Value = 17
$ws.Range("AG18").Value = 301
$ws.Range("AI18").Value = 21
$ws.Range("AN18").Value = 3.75
$ws.Range("AX18").Value = 23
$ws.Range("G18").Value = 1.85
$ws.Range("H18").Value = 3.75
$ws.Range("I18").Value = 4.2
$ws.Range("L18").Value = 4.75
$ws.Range("Q18").Value = 2
$ws.Range("R18").Value = 1.9
$ws.Range("X18").Value = 8.5
$ws.Range("Q19").Value = 2.15
$ws.Range("R19").Value = 1.67
$ws.Range("BB20").Value = 151
$ws.Range("G20").Value = 2.9
$ws.Range("I20").Value = 2.55
$ws.Range("Q21").Value = 2.1
$ws.Range("R21").Value = 1.7
$ws.Range("AC22").Value = 10
$ws.Range("AD22").Value = 6
$ws.Range("AU22").Value = 8
$ws.Range("AW22").Value = 4.5
$ws.Range("AZ22").Value = 51
$ws.Range("G22").Value = 2.75
$ws.Range("I22").Value = 2.55
$ws.Range("J22").Value = 3.4
$ws.Range("K22").Value = 2.1
$ws.Range("M22").Value = 1.06
$ws.Range("N22").Value = 10
$ws.Range("O22").Value = 1.3
$ws.Range("P22").Value = 3.5
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = 1.9
$ws.Range("U22").Value = 1.75
$ws.Range("V22").Value = 2
$ws.Range("Z22").Value = 29
